$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Fstl1"
$ws.Cells.Item(2, 3).Value = "Dip2a"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 15.09508766666667
$ws.Cells.Item(2, 8).Value = 45.285263
$ws.Cells.Item(2, 9).Value = 0.03390535125908365
$ws.Cells.Item(2, 10).Value = 0.03390535125908365
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 7.306170666666667
$ws.Cells.Item(2, 14).Value = 21.918512
$ws.Cells.Item(2, 15).Value = 0.3391103748162015
$ws.Cells.Item(2, 16).Value = 0.3391103748162015
$ws.Cells.Item(2, 17).Value = 110.2872867209618
$ws.Cells.Item(2, 18).Value = 992.585580488656
$ws.Cells.Item(2, 19).Value = 0.01149765637374283
$ws.Cells.Item(2, 20).Value = 0.01149765637374283

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Fstl1"
$ws.Cells.Item(3, 3).Value = "Dip2a"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 15.09508766666667
$ws.Cells.Item(3, 8).Value = 45.285263
$ws.Cells.Item(3, 9).Value = 0.03390535125908365
$ws.Cells.Item(3, 10).Value = 0.03390535125908365
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 9.033654666666669
$ws.Cells.Item(3, 14).Value = 27.100964
$ws.Cells.Item(3, 15).Value = 0.4192902355744033
$ws.Cells.Item(3, 16).Value = 0.4192902355744033
$ws.Cells.Item(3, 17).Value = 136.3638091437258
$ws.Cells.Item(3, 18).Value = 1227.274282293532
$ws.Cells.Item(3, 19).Value = 0.01421618271665408
$ws.Cells.Item(3, 20).Value = 0.01421618271665408

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Fstl1"
$ws.Cells.Item(4, 3).Value = "Dip2a"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 15.09508766666667
$ws.Cells.Item(4, 8).Value = 45.285263
$ws.Cells.Item(4, 9).Value = 0.03390535125908365
$ws.Cells.Item(4, 10).Value = 0.03390535125908365
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 5.205285666666666
$ws.Cells.Item(4, 14).Value = 15.615857
$ws.Cells.Item(4, 15).Value = 0.2415993896093951
$ws.Cells.Item(4, 16).Value = 0.2415993896093952
$ws.Cells.Item(4, 17).Value = 78.57424346837676
$ws.Cells.Item(4, 18).Value = 707.168191215391
$ws.Cells.Item(4, 19).Value = 0.008191512168686747
$ws.Cells.Item(4, 20).Value = 0.008191512168686748

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Fstl1"
$ws.Cells.Item(5, 3).Value = "Dip2a"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 398.9908546666666
$ws.Cells.Item(5, 8).Value = 1196.972564
$ws.Cells.Item(5, 9).Value = 0.8961806234824337
$ws.Cells.Item(5, 10).Value = 0.8961806234824337
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 7.306170666666667
$ws.Cells.Item(5, 14).Value = 21.918512
$ws.Cells.Item(5, 15).Value = 0.3391103748162015
$ws.Cells.Item(5, 16).Value = 0.3391103748162015
$ws.Cells.Item(5, 17).Value = 2915.095278633863
$ws.Cells.Item(5, 18).Value = 26235.85750770476
$ws.Cells.Item(5, 19).Value = 0.3039041471321452
$ws.Cells.Item(5, 20).Value = 0.3039041471321453

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Fstl1"
$ws.Cells.Item(6, 3).Value = "Dip2a"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 398.9908546666666
$ws.Cells.Item(6, 8).Value = 1196.972564
$ws.Cells.Item(6, 9).Value = 0.8961806234824337
$ws.Cells.Item(6, 10).Value = 0.8961806234824337
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 9.033654666666669
$ws.Cells.Item(6, 14).Value = 27.100964
$ws.Cells.Item(6, 15).Value = 0.4192902355744033
$ws.Cells.Item(6, 16).Value = 0.4192902355744033
$ws.Cells.Item(6, 17).Value = 3604.345596216855
$ws.Cells.Item(6, 18).Value = 32439.1103659517
$ws.Cells.Item(6, 19).Value = 0.3757597847371653
$ws.Cells.Item(6, 20).Value = 0.3757597847371653

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Fstl1"
$ws.Cells.Item(7, 3).Value = "Dip2a"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 398.9908546666666
$ws.Cells.Item(7, 8).Value = 1196.972564
$ws.Cells.Item(7, 9).Value = 0.8961806234824337
$ws.Cells.Item(7, 10).Value = 0.8961806234824337
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 5.205285666666666
$ws.Cells.Item(7, 14).Value = 15.615857
$ws.Cells.Item(7, 15).Value = 0.2415993896093951
$ws.Cells.Item(7, 16).Value = 0.2415993896093952
$ws.Cells.Item(7, 17).Value = 2076.861376927483
$ws.Cells.Item(7, 18).Value = 18691.75239234735
$ws.Cells.Item(7, 19).Value = 0.2165166916131231
$ws.Cells.Item(7, 20).Value = 0.2165166916131232

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Fstl1"
$ws.Cells.Item(8, 3).Value = "Dip2a"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 31.12660099999999
$ws.Cells.Item(8, 8).Value = 93.37980299999998
$ws.Cells.Item(8, 9).Value = 0.06991402525848273
$ws.Cells.Item(8, 10).Value = 0.06991402525848271
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 7.306170666666667
$ws.Cells.Item(8, 14).Value = 21.918512
$ws.Cells.Item(8, 15).Value = 0.3391103748162015
$ws.Cells.Item(8, 16).Value = 0.3391103748162015
$ws.Cells.Item(8, 17).Value = 227.4162591792373
$ws.Cells.Item(8, 18).Value = 2046.746332613136
$ws.Cells.Item(8, 19).Value = 0.02370857131031346
$ws.Cells.Item(8, 20).Value = 0.02370857131031345

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Fstl1"
$ws.Cells.Item(9, 3).Value = "Dip2a"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 31.12660099999999
$ws.Cells.Item(9, 8).Value = 93.37980299999998
$ws.Cells.Item(9, 9).Value = 0.06991402525848273
$ws.Cells.Item(9, 10).Value = 0.06991402525848271
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 9.033654666666669
$ws.Cells.Item(9, 14).Value = 27.100964
$ws.Cells.Item(9, 15).Value = 0.4192902355744033
$ws.Cells.Item(9, 16).Value = 0.4192902355744033
$ws.Cells.Item(9, 17).Value = 281.1869643811214
$ws.Cells.Item(9, 18).Value = 2530.682679430092
$ws.Cells.Item(9, 19).Value = 0.02931426812058401
$ws.Cells.Item(9, 20).Value = 0.029314268120584

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Fstl1"
$ws.Cells.Item(10, 3).Value = "Dip2a"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 31.12660099999999
$ws.Cells.Item(10, 8).Value = 93.37980299999998
$ws.Cells.Item(10, 9).Value = 0.06991402525848273
$ws.Cells.Item(10, 10).Value = 0.06991402525848271
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 5.205285666666666
$ws.Cells.Item(10, 14).Value = 15.615857
$ws.Cells.Item(10, 15).Value = 0.2415993896093951
$ws.Cells.Item(10, 16).Value = 0.2415993896093952
$ws.Cells.Item(10, 17).Value = 162.0228500373523
$ws.Cells.Item(10, 18).Value = 1458.20565033617
$ws.Cells.Item(10, 19).Value = 0.01689118582758526
$ws.Cells.Item(10, 20).Value = 0.01689118582758526
